# Apply swapped/updated field data for the Artfynd observation rows.
# The underlying export re-ordered/reassigned several observation records
# (same dataset, rows shuffled + a few field corrections), so each affected
# row is rewritten cell-by-cell with its new content.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 11 ---
$ws.Range("A11").Value = 131085696
$ws.Range("B11").Value = 57884
$ws.Range("E11").Value = 100109
$ws.Range("F11").Value = 'Tretåig hackspett'
$ws.Range("G11").Value = 'Picoides tridactylus'
$ws.Range("H11").Value = '(Linnaeus, 1758)'
$ws.Range("M11").Value = 'färska spår'
$ws.Range("Q11").Value = 585207
$ws.Range("R11").Value = 7060471
$ws.Range("S11").Value = 10
$ws.Range("Z11").Value = '11:55'
$ws.Range("AB11").Value = '11:55'
$ws.Range("AC11").Value = 'Ringhack på gran'
$ws.Range("AW11").Value = 'Kim Hultgren'
$ws.Range("AX11").Value = 'Kim Hultgren'

# --- Row 12 ---
$ws.Range("A12").Value = 131087481
$ws.Range("B12").Value = 91828
$ws.Range("E12").Value = 5432
$ws.Range("F12").Value = 'Granticka'
$ws.Range("G12").Value = 'Porodaedalea chrysoloma s.lat.'
$ws.Range("H12").Value = ""
$ws.Range("M12").Value = ""
$ws.Range("Q12").Value = 585150
$ws.Range("R12").Value = 7060657
$ws.Range("S12").Value = 15
$ws.Range("Z12").Value = ""
$ws.Range("AB12").Value = ""
$ws.Range("AC12").Value = ""
$ws.Range("AW12").Value = 'Daniel Rutschman'
$ws.Range("AX12").Value = 'Daniel Rutschman'

# --- Row 18 ---
$ws.Range("A18").Value = 131092560
$ws.Range("B18").Value = 91804
$ws.Range("E18").Value = 1108
$ws.Range("F18").Value = 'Harticka'
$ws.Range("G18").Value = 'Pelloporus leporinus'
$ws.Range("H18").Value = '(Fr.) Krieglst.'
$ws.Range("M18").Value = ""
$ws.Range("Q18").Value = 585129
$ws.Range("R18").Value = 7060254
$ws.Range("S18").Value = 10
$ws.Range("Z18").Value = '15:17'
$ws.Range("AB18").Value = '15:17'
$ws.Range("AC18").Value = ""
$ws.Range("AW18").Value = 'Kim Hultgren'
$ws.Range("AX18").Value = 'Kim Hultgren'

# --- Row 19 ---
$ws.Range("A19").Value = 131092554
$ws.Range("B19").Value = 57884
$ws.Range("E19").Value = 100109
$ws.Range("F19").Value = 'Tretåig hackspett'
$ws.Range("G19").Value = 'Picoides tridactylus'
$ws.Range("H19").Value = '(Linnaeus, 1758)'
$ws.Range("M19").Value = 'äldre spår'
$ws.Range("Q19").Value = 585147
$ws.Range("R19").Value = 7060312
$ws.Range("S19").Value = 15
$ws.Range("Z19").Value = ""
$ws.Range("AB19").Value = ""
$ws.Range("AC19").Value = 'Äldre ringhack, gran'
$ws.Range("AW19").Value = 'Daniel Rutschman'
$ws.Range("AX19").Value = 'Daniel Rutschman'

# --- Row 20 ---
$ws.Range("A20").Value = 131086957
$ws.Range("B20").Value = 57884
$ws.Range("E20").Value = 100109
$ws.Range("F20").Value = 'Tretåig hackspett'
$ws.Range("G20").Value = 'Picoides tridactylus'
$ws.Range("H20").Value = '(Linnaeus, 1758)'
$ws.Range("M20").Value = 'färska spår'
$ws.Range("Q20").Value = 585162
$ws.Range("R20").Value = 7060573
$ws.Range("S20").Value = 10
$ws.Range("Z20").Value = '12:21'
$ws.Range("AB20").Value = '12:21'
$ws.Range("AC20").Value = 'Ringhack på gran'
$ws.Range("AW20").Value = 'Kim Hultgren'
$ws.Range("AX20").Value = 'Kim Hultgren'

# --- Row 21 ---
$ws.Range("A21").Value = 131085737
$ws.Range("B21").Value = 79243
$ws.Range("E21").Value = 6425
$ws.Range("F21").Value = 'Garnlav'
$ws.Range("G21").Value = 'Alectoria sarmentosa'
$ws.Range("H21").Value = '(Ach.) Ach.'
$ws.Range("M21").Value = ""
$ws.Range("Q21").Value = 585170
$ws.Range("R21").Value = 7060469
$ws.Range("S21").Value = 15
$ws.Range("Z21").Value = '11:58'
$ws.Range("AB21").Value = '11:58'
$ws.Range("AC21").Value = ""
$ws.Range("AW21").Value = 'Daniel Rutschman'
$ws.Range("AX21").Value = 'Daniel Rutschman'

# --- Row 22 ---
$ws.Range("A22").Value = 131085446
$ws.Range("B22").Value = 79243
$ws.Range("E22").Value = 6425
$ws.Range("F22").Value = 'Garnlav'
$ws.Range("G22").Value = 'Alectoria sarmentosa'
$ws.Range("H22").Value = '(Ach.) Ach.'
$ws.Range("M22").Value = ""
$ws.Range("Q22").Value = 585301
$ws.Range("R22").Value = 7060488
$ws.Range("S22").Value = 10
$ws.Range("Z22").Value = '11:41'
$ws.Range("AB22").Value = '11:41'
$ws.Range("AC22").Value = ""
$ws.Range("AW22").Value = 'Kim Hultgren'
$ws.Range("AX22").Value = 'Kim Hultgren'

# --- Row 28 ---
$ws.Range("A28").Value = 131085484
$ws.Range("B28").Value = 57884
$ws.Range("E28").Value = 100109
$ws.Range("F28").Value = 'Tretåig hackspett'
$ws.Range("G28").Value = 'Picoides tridactylus'
$ws.Range("H28").Value = '(Linnaeus, 1758)'
$ws.Range("M28").Value = 'färska spår'
$ws.Range("Q28").Value = 585303
$ws.Range("R28").Value = 7060488
$ws.Range("S28").Value = 15
$ws.Range("Z28").Value = ""
$ws.Range("AB28").Value = ""
$ws.Range("AC28").Value = 'Färska ringhack gran'
$ws.Range("AW28").Value = 'Daniel Rutschman'
$ws.Range("AX28").Value = 'Daniel Rutschman'

# --- Row 29 ---
$ws.Range("A29").Value = 131085240
$ws.Range("B29").Value = 57884
$ws.Range("E29").Value = 100109
$ws.Range("F29").Value = 'Tretåig hackspett'
$ws.Range("G29").Value = 'Picoides tridactylus'
$ws.Range("H29").Value = '(Linnaeus, 1758)'
$ws.Range("M29").Value = 'färska spår'
$ws.Range("Q29").Value = 585289
$ws.Range("R29").Value = 7060293
$ws.Range("S29").Value = 10
$ws.Range("Z29").Value = '11:16'
$ws.Range("AB29").Value = '11:16'
$ws.Range("AC29").Value = 'Ringhack på tall'
$ws.Range("AW29").Value = 'Kim Hultgren'
$ws.Range("AX29").Value = 'Kim Hultgren'

# --- Row 30 ---
$ws.Range("A30").Value = 131085171
$ws.Range("B30").Value = 91804
$ws.Range("E30").Value = 1108
$ws.Range("F30").Value = 'Harticka'
$ws.Range("G30").Value = 'Pelloporus leporinus'
$ws.Range("H30").Value = '(Fr.) Krieglst.'
$ws.Range("M30").Value = ""
$ws.Range("Q30").Value = 585222
$ws.Range("R30").Value = 7060254
$ws.Range("S30").Value = 15
$ws.Range("Z30").Value = ""
$ws.Range("AB30").Value = ""
$ws.Range("AC30").Value = ""
$ws.Range("AW30").Value = 'Daniel Rutschman'
$ws.Range("AX30").Value = 'Daniel Rutschman'

# --- Row 31 ---
$ws.Range("A31").Value = 131085178
$ws.Range("B31").Value = 91828
$ws.Range("E31").Value = 5432
$ws.Range("F31").Value = 'Granticka'
$ws.Range("G31").Value = 'Porodaedalea chrysoloma s.lat.'
$ws.Range("H31").Value = ""
$ws.Range("M31").Value = ""
$ws.Range("Q31").Value = 585225
$ws.Range("R31").Value = 7060258
$ws.Range("S31").Value = 10
$ws.Range("Z31").Value = '11:08'
$ws.Range("AB31").Value = '11:08'
$ws.Range("AC31").Value = ""
$ws.Range("AW31").Value = 'Kim Hultgren'
$ws.Range("AX31").Value = 'Kim Hultgren'

# --- Row 32 ---
$ws.Range("A32").Value = 131085569
$ws.Range("B32").Value = 79243
$ws.Range("E32").Value = 6425
$ws.Range("F32").Value = 'Garnlav'
$ws.Range("G32").Value = 'Alectoria sarmentosa'
$ws.Range("H32").Value = '(Ach.) Ach.'
$ws.Range("M32").Value = ""
$ws.Range("Q32").Value = 585249
$ws.Range("R32").Value = 7060505
$ws.Range("S32").Value = 15
$ws.Range("Z32").Value = ""
$ws.Range("AB32").Value = ""
$ws.Range("AC32").Value = ""
$ws.Range("AW32").Value = 'Daniel Rutschman'
$ws.Range("AX32").Value = 'Daniel Rutschman'

# --- Row 33 ---
$ws.Range("A33").Value = 131087388
$ws.Range("B33").Value = 79243
$ws.Range("E33").Value = 6425
$ws.Range("F33").Value = 'Garnlav'
$ws.Range("G33").Value = 'Alectoria sarmentosa'
$ws.Range("H33").Value = '(Ach.) Ach.'
$ws.Range("M33").Value = ""
$ws.Range("Q33").Value = 585131
$ws.Range("R33").Value = 7060627
$ws.Range("S33").Value = 15
$ws.Range("Z33").Value = ""
$ws.Range("AB33").Value = ""
$ws.Range("AC33").Value = ""
$ws.Range("AW33").Value = 'Daniel Rutschman'
$ws.Range("AX33").Value = 'Daniel Rutschman'

# --- Row 34 ---
$ws.Range("A34").Value = 131092590
$ws.Range("B34").Value = 79243
$ws.Range("E34").Value = 6425
$ws.Range("F34").Value = 'Garnlav'
$ws.Range("G34").Value = 'Alectoria sarmentosa'
$ws.Range("H34").Value = '(Ach.) Ach.'
$ws.Range("M34").Value = ""
$ws.Range("Q34").Value = 585145
$ws.Range("R34").Value = 7060230
$ws.Range("S34").Value = 10
$ws.Range("Z34").Value = '15:20'
$ws.Range("AB34").Value = '15:20'
$ws.Range("AC34").Value = ""
$ws.Range("AW34").Value = 'Kim Hultgren'
$ws.Range("AX34").Value = 'Kim Hultgren'

# --- Row 35 ---
$ws.Range("A35").Value = 131085126
$ws.Range("B35").Value = 57064
$ws.Range("E35").Value = 102612
$ws.Range("F35").Value = 'Järpe'
$ws.Range("G35").Value = 'Tetrastes bonasia'
$ws.Range("H35").Value = '(Linnaeus, 1758)'
$ws.Range("M35").Value = 'färsk spillning'
$ws.Range("Q35").Value = 585219
$ws.Range("R35").Value = 7060240
$ws.Range("S35").Value = 15
$ws.Range("Z35").Value = ""
$ws.Range("AB35").Value = ""
$ws.Range("AC35").Value = ""
$ws.Range("AW35").Value = 'Daniel Rutschman'
$ws.Range("AX35").Value = 'Daniel Rutschman'

# --- Row 36 ---
$ws.Range("A36").Value = 131092585
$ws.Range("B36").Value = 91804
$ws.Range("E36").Value = 1108
$ws.Range("F36").Value = 'Harticka'
$ws.Range("G36").Value = 'Pelloporus leporinus'
$ws.Range("H36").Value = '(Fr.) Krieglst.'
$ws.Range("M36").Value = ""
$ws.Range("Q36").Value = 585130
$ws.Range("R36").Value = 7060263
$ws.Range("S36").Value = 15
$ws.Range("Z36").Value = ""
$ws.Range("AB36").Value = ""
$ws.Range("AC36").Value = ""
$ws.Range("AW36").Value = 'Daniel Rutschman'
$ws.Range("AX36").Value = 'Daniel Rutschman'

